# "merge with random planar graph"
# Replace the i/j/num edge-list data below the header row with the new,
# larger merged edge list (now 54 data rows instead of 40).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("G","P",20),
    @("D","E",20),
    @("M","T",20),
    @("F","H",20),
    @("D","Q",20),
    @("E","S",20),
    @("A","F",20),
    @("J","L",20),
    @("I","P",20),
    @("A","L",20),
    @("K","P",20),
    @("C","R",20),
    @("B","P",20),
    @("H","Q",20),
    @("G","O",20),
    @("E","O",20),
    @("D","G",20),
    @("C","E",20),
    @("B","C",20),
    @("A","E",20),
    @("L","N",20),
    @("L","T",20),
    @("B","R",20),
    @("J","M",20),
    @("B","K",20),
    @("C","S",20),
    @("G","P",20),
    @("D","E",20),
    @("M","T",20),
    @("F","H",20),
    @("D","Q",20),
    @("E","S",20),
    @("A","F",20),
    @("J","L",20),
    @("I","P",20),
    @("A","L",20),
    @("K","P",20),
    @("C","R",20),
    @("B","P",20),
    @("H","Q",20),
    @("G","O",20),
    @("E","O",20),
    @("D","G",20),
    @("A","H",20),
    @("C","E",20),
    @("B","C",20),
    @("A","E",20),
    @("L","N",20),
    @("L","T",20),
    @("B","I",20),
    @("B","R",20),
    @("J","M",20),
    @("B","K",20),
    @("C","S",20)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
